$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Calculated CPI" header to "CPI using Miss Rate"
$ws.Range("H1").Value = "CPI using Miss Rate"

# Add new header "CPI using Miss Ticks" in column I, copying the bold/
# bordered header style used by the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "CPI using Miss Ticks"

# New CPI (using miss ticks) values for rows 2-87 (column I)
$missTicksCpi = @(
    1.019088819932465,
    1.023137235539063,
    1.024614725353648,
    1.024654752962982,
    1.024672631043375,
    1.024682641962886,
    1.024676198087543,
    1.024689347653919,
    1.213642583497051,
    1.094881339048007,
    1.046702497419491,
    1.024654752962982,
    1.01817412638853,
    1.010547152178551,
    1.008704953243769,
    1.005610650822219,
    1.024710545718464,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.083759814204811,
    1.052942910123217,
    1.031939664294966,
    1.024654752962982,
    1.026492576786782,
    1.037541853663117,
    1.07572697587002,
    1.024610014196572,
    1.024654752962982,
    1.024662076996317,
    1.024636785081533,
    1.137561197819615,
    1.024654752962982,
    1.012298486563128,
    1.011800945835822,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.024654752962982,
    1.818778183546783,
    1.818914285413227,
    1.818924458847237,
    1.818893596775429,
    1.818884338046038,
    1.818846553453971,
    1.819041641142949,
    1.819049030194306,
    1.819736066944585,
    1.819053449046702,
    1.819066683280637,
    1.818893596775429,
    1.818851959404859,
    1.819156605316266,
    1.819156556194362,
    1.819156538541178,
    1.819178433184797,
    1.818831369699976,
    1.819225379348832,
    1.818893596775429,
    1.818672046427031,
    1.818226414992708,
    1.817404032428819,
    1.788163666731069,
    1.956703954345838,
    1.953237169491113,
    1.902151521661942,
    1.818893596775429,
    1.713588805558187,
    1.619220111036396,
    1.48165700843835,
    1.818813880125368,
    1.818893596775429,
    1.818978096568205,
    1.818998541685307,
    1.819064990520892,
    1.818893596775429,
    1.818869442566909,
    1.818869755447212,
    1.818722761707717,
    1.818893596775429,
    1.819168757390788,
    1.818951804196615
)

for ($i = 0; $i -lt $missTicksCpi.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $missTicksCpi[$i]
}

Write-Host "Added CPI using Miss Ticks column (I) with $($missTicksCpi.Length) data rows"
